# Generate Report for Handoff
# Adds the newly-handed-off file "e07622b0-cb26-40e5-98df-890a223460f3.md"
# as a new row (row 9) to all three tables: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$fileBase    = "e07622b0-cb26-40e5-98df-890a223460f3"
$newMdFile   = "$fileBase.md"
$commitSha   = "6d907d3ea5a81a31ab1f800f4214f1fcf8233c76"
$srcUrl      = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$commitSha/e2e/$newMdFile"

# ---------------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.ListRows.Add() | Out-Null

$wsOverview.Range("A9").Value = $newMdFile
$wsOverview.Hyperlinks.Add($wsOverview.Range("B9"), $srcUrl, "", "", "e2e\$newMdFile") | Out-Null
$wsOverview.Range("C9").Value = "'.md"
$wsOverview.Range("D9").Value = "'"
$wsOverview.Range("E9").Value = "'Ready for handoff"
$wsOverview.Range("F9").Value = "'Ready for handoff"
$wsOverview.Range("G9").Value = "'2016-09-01 16:52:39"

# ---------------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.ListRows.Add() | Out-Null

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A9"), $srcUrl, "", "", $newMdFile) | Out-Null
$wsZhCn.Range("B9").Value = "'.md"
$wsZhCn.Range("C9").Value = "'Ready for handoff"
$wsZhCn.Range("D9").Value = "'e2e"
$wsZhCn.Range("E9").Value = "'ht"
$wsZhCn.Range("F9").Value = "'False"
$wsZhCn.Range("G9").Value = "'$fileBase.$commitSha.zh-cn.xlf"
$wsZhCn.Range("H9").Value = "'2016-09-01 16:52:34"
$wsZhCn.Range("I9").Value = "'"
$wsZhCn.Range("J9").Value = "'"
$wsZhCn.Range("K9").Value = "'0001-01-01 00:00:00"
$wsZhCn.Range("L9").Value = "'"
$wsZhCn.Range("M9").Value = "'True"
$wsZhCn.Range("N9").Value = "'"
$wsZhCn.Range("O9").Value = "'False"
$wsZhCn.Range("P9").Value = "'"

# ---------------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.ListRows.Add() | Out-Null

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A9"), $srcUrl, "", "", $newMdFile) | Out-Null
$wsDeDe.Range("B9").Value = "'.md"
$wsDeDe.Range("C9").Value = "'Ready for handoff"
$wsDeDe.Range("D9").Value = "'e2e"
$wsDeDe.Range("E9").Value = "'ht"
$wsDeDe.Range("F9").Value = "'False"
$wsDeDe.Range("G9").Value = "'$fileBase.$commitSha.de-de.xlf"
$wsDeDe.Range("H9").Value = "'2016-09-01 16:52:39"
$wsDeDe.Range("I9").Value = "'"
$wsDeDe.Range("J9").Value = "'"
$wsDeDe.Range("K9").Value = "'0001-01-01 00:00:00"
$wsDeDe.Range("L9").Value = "'"
$wsDeDe.Range("M9").Value = "'True"
$wsDeDe.Range("N9").Value = "'"
$wsDeDe.Range("O9").Value = "'False"
$wsDeDe.Range("P9").Value = "'"

Write-Host "Added handoff row for $newMdFile to Overview, zh-cn and de-de sheets."
